{"js": "// Remove the table row describing \"O Sistema deve ser acess\u00edvel de m\u00faltiplas\n// plataformas\" / \"Necessidade dos Stakeholders...\" and merge the two\n// adjacent runs that make up \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de\n// licen\u00e7as de uso ou royalties.\" into a single run.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Locate the row to delete by its first-cell text, rather than assuming a\n// fixed index, so the script is resilient to table shape.\nlet rowToDelete = null;\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  const firstCell = row.cells.items[0];\n  firstCell.body.load(\"text\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  const firstCell = row.cells.items[0];\n  if (firstCell.body.text.trim() === \"O Sistema deve ser acess\u00edvel de m\u00faltiplas plataformas\") {\n    rowToDelete = row;\n    break;\n  }\n}\n\nif (rowToDelete) {\n  rowToDelete.delete();\n  await context.sync();\n}\n\n// Merge the two runs \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de licen\u00e7as de\n// uso ou royalties\" + \".\" into a single run by replacing the matched range's\n// text in place (this collapses the split runs into one run using the\n// formatting of the search hit).\nconst searchResults = context.document.body.search(\n  \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de licen\u00e7as de uso ou royalties.\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const target = searchResults.items[0];\n  target.insertText(\n    \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de licen\u00e7as de uso ou royalties.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Remove the table row describing \"O Sistema deve ser acess\u00edvel de m\u00faltiplas\n# plataformas\" / \"Necessidade dos Stakeholders...\" and merge the two\n# adjacent runs that make up \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de\n# licen\u00e7as de uso ou royalties.\" into a single run.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Find the row whose first cell holds the restriction to remove, rather than\n# assuming a fixed row index.\n$targetText = \"O Sistema deve ser acess\u00edvel de m\u00faltiplas plataformas\"\nfor ($i = $table.Rows.Count; $i -ge 1; $i--) {\n    $cellText = $table.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($cellText -eq $targetText) {\n        $table.Rows.Item($i).Delete()\n        break\n    }\n}\n\n# Merge the two runs \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de licen\u00e7as de\n# uso ou royalties\" + \".\" into a single run via Find & Replace, which\n# collapses the matched range into one run using the formatting already\n# present there.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$searchText = \"cuja utiliza\u00e7\u00e3o n\u00e3o implica o pagamento de licen\u00e7as de uso ou royalties.\"\n$find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null\n"}
